# Auto-generated Excel COM-interop script applying numeric corrections
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 41673076
$ws.Range("I69").Value = 166669920
$ws.Range("J69").Value = 7458.3335
$ws.Range("K69").Value = 500009760
$ws.Range("L69").Value = 22375.0005
$ws.Range("M69").Value = -500008886
$ws.Range("N69").Value = -24123.0005
$ws.Range("H72").Value = 41673076
$ws.Range("I72").Value = 166669920
$ws.Range("J72").Value = 7458.3335
$ws.Range("K72").Value = 1500029280
$ws.Range("L72").Value = 67125.0015
$ws.Range("M72").Value = -1500024912
$ws.Range("N72").Value = -75861.0015
$ws.Range("H106").Value = 3272.3333
$ws.Range("I106").Value = 3048.9092
$ws.Range("K106").Value = 3048.9092
$ws.Range("M106").Value = -2417.9092
$ws.Range("H113").Value = 43755.707
$ws.Range("I113").Value = 2674.9
$ws.Range("J113").Value = 102442.57
$ws.Range("K113").Value = 2674.9
$ws.Range("L113").Value = 102442.57
$ws.Range("M113").Value = 579.0999999999999
$ws.Range("N113").Value = -108950.57

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2916.33
$ws.Range("I32").Value = 2508.5103
$ws.Range("K32").Value = 2508.5103
$ws.Range("M32").Value = -2221.5103
$ws.Range("H61").Value = 7321.95
$ws.Range("I61").Value = 7246.6113
$ws.Range("K61").Value = 7246.6113
$ws.Range("M61").Value = -7034.6113
$ws.Range("H97").Value = 647.2174
$ws.Range("I97").Value = 461
$ws.Range("J97").Value = 1174.8334
$ws.Range("K97").Value = 461
$ws.Range("L97").Value = 1174.8334
$ws.Range("M97").Value = 35
$ws.Range("N97").Value = -2166.8334
$ws.Range("H132").Value = 8513.558000000001
$ws.Range("I132").Value = 8135.9146
$ws.Range("K132").Value = 24407.7438
$ws.Range("M132").Value = -21877.7438
$ws.Range("H136").Value = 7321.95
$ws.Range("I136").Value = 7246.6113
$ws.Range("K136").Value = 21739.8339
$ws.Range("M136").Value = -19189.8339
$ws.Range("H139").Value = 117142.336
$ws.Range("J139").Value = 117142.336
$ws.Range("L139").Value = 117142.336
$ws.Range("N139").Value = -127422.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5694.9443
$ws.Range("I20").Value = 8009.5
$ws.Range("K20").Value = 8009.5
$ws.Range("M20").Value = -7762.5
$ws.Range("H86").Value = 3369.138
$ws.Range("I86").Value = 3093.875
$ws.Range("K86").Value = 3093.875
$ws.Range("M86").Value = -1970.875
$ws.Range("H89").Value = 3369.138
$ws.Range("I89").Value = 3093.875
$ws.Range("K89").Value = 15469.375
$ws.Range("M89").Value = -9853.375
$ws.Range("H105").Value = 3657.7693
$ws.Range("I105").Value = 3713.9092
$ws.Range("K105").Value = 3713.9092
$ws.Range("M105").Value = -1966.9092
$ws.Range("H107").Value = 2295.7073
$ws.Range("I107").Value = 1784.3429
$ws.Range("K107").Value = 1784.3429
$ws.Range("M107").Value = 135.6570999999999
$ws.Range("H134").Value = 3125.0938
$ws.Range("I134").Value = 3000.8
$ws.Range("K134").Value = 9002.400000000001
$ws.Range("M134").Value = -6467.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2339.7585
$ws.Range("I58").Value = 2336.6667
$ws.Range("J58").Value = 2344.818
$ws.Range("K58").Value = 2336.6667
$ws.Range("L58").Value = 2344.818
$ws.Range("M58").Value = -2133.6667
$ws.Range("N58").Value = -2750.818
$ws.Range("H122").Value = 2434
$ws.Range("I122").Value = 2016.4
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 6049.200000000001
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -3599.200000000001
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 38097732
$ws.Range("I132").Value = 45979040
$ws.Range("J132").Value = 4731.1665
$ws.Range("K132").Value = 137937120
$ws.Range("L132").Value = 14193.4995
$ws.Range("M132").Value = -137934590
$ws.Range("N132").Value = -19253.4995
$ws.Range("H136").Value = 2339.7585
$ws.Range("I136").Value = 2336.6667
$ws.Range("J136").Value = 2344.818
$ws.Range("K136").Value = 7010.000100000001
$ws.Range("L136").Value = 7034.454000000001
$ws.Range("M136").Value = -4460.000100000001
$ws.Range("N136").Value = -12134.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 316.72223
$ws.Range("I14").Value = 316.72223
$ws.Range("K14").Value = 950.16669
$ws.Range("M14").Value = -777.16669
$ws.Range("H57").Value = 9749.333000000001
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 9749.333000000001
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 29247.999
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -30365.999
$ws.Range("H107").Value = 752.6667
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 784.25
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 2352.75
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -6192.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1231.3334
$ws.Range("I97").Value = 981.88464
$ws.Range("J97").Value = 2852.75
$ws.Range("K97").Value = 981.88464
$ws.Range("L97").Value = 2852.75
$ws.Range("M97").Value = -485.88464
$ws.Range("N97").Value = -3844.75
$ws.Range("H102").Value = 11596149
$ws.Range("I102").Value = 12755289
$ws.Range("K102").Value = 12755289
$ws.Range("M102").Value = -12753667
$ws.Range("H113").Value = 1032.381
$ws.Range("I113").Value = 1040.5
$ws.Range("K113").Value = 1040.5
$ws.Range("M113").Value = 1129.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1652.5151
$ws.Range("J16").Value = 3933.1667
$ws.Range("L16").Value = 3933.1667
$ws.Range("N16").Value = -4273.1667
$ws.Range("H82").Value = 2241.2144
$ws.Range("I82").Value = 1163.6875
$ws.Range("K82").Value = 1163.6875
$ws.Range("M82").Value = -802.6875
$ws.Range("H85").Value = 2241.2144
$ws.Range("I85").Value = 1163.6875
$ws.Range("K85").Value = 1163.6875
$ws.Range("M85").Value = 84.3125
$ws.Range("H93").Value = 1500.4445
$ws.Range("I93").Value = 1317.7222
$ws.Range("K93").Value = 1317.7222
$ws.Range("M93").Value = -69.72219999999993
$ws.Range("H136").Value = 4278.3335
$ws.Range("I136").Value = 2667.625
$ws.Range("K136").Value = 8002.875
$ws.Range("M136").Value = -5452.875
$ws.Range("H140").Value = 58439.07
$ws.Range("J140").Value = 58439.07
$ws.Range("L140").Value = 58439.07
$ws.Range("N140").Value = -68799.07000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 46731.473
$ws.Range("I2").Value = 46731.473
$ws.Range("K2").Value = 46731.473
$ws.Range("M2").Value = -46619.473
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 53998.5
$ws.Range("J49").Value = 53998.5
$ws.Range("L49").Value = 53998.5
$ws.Range("N49").Value = -54458.5
$ws.Range("H54").Value = 38594.6
$ws.Range("J54").Value = 38243.25
$ws.Range("L54").Value = 38243.25
$ws.Range("N54").Value = -39283.25
$ws.Range("H106").Value = 120000
$ws.Range("J106").Value = 120000
$ws.Range("L106").Value = 120000
$ws.Range("N106").Value = -122524
$ws.Range("I132").Value = 13890332
$ws.Range("K132").Value = 41670996
$ws.Range("M132").Value = -41668466
$ws.Range("H136").Value = 6116.278
$ws.Range("I136").Value = 4461.905
$ws.Range("J136").Value = 8432.4
$ws.Range("K136").Value = 13385.715
$ws.Range("L136").Value = 25297.2
$ws.Range("M136").Value = -10835.715
$ws.Range("N136").Value = -30397.2
